$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Split the sentence "...deseja visualizar algum usuário, cadastrar,
#    editar ou remover o mesmo." into "...deseja" + " procurar por um
#    usuário especifico." as two separate runs with identical formatting
#    (matching the run-per-run structure produced by the original edit).
# -----------------------------------------------------------------------
$oldSentence = "ema deseja visualizar algum usuário, cadastrar, editar ou remover o mesmo."

$find1 = $d.Content
$found1 = $find1.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false)
if (-not $found1) {
    throw "Could not find the target sentence to edit."
}
$matchStart = $find1.Start
$matchEnd = $find1.End

# Use InsertXML on a *freshly constructed* Range (the Range handed back
# by Find must not be reused directly - it behaves as an insertion point
# rather than a replacement target) so we can precisely control the
# resulting run split, including preserving the original run's rsid.
$replaceRange = $d.Range($matchStart, $matchEnd)
$runsXml = '<w:r w:rsidR="00BE0992"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="ar-SA"/></w:rPr><w:t>ema deseja</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> procurar por um usuário especifico.</w:t></w:r>'
$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$replaceRange.InsertXML($packageXml)

# -----------------------------------------------------------------------
# 2) "Pontos de Inclusão" bullet list:
#      Gerar Relatórios   -> Cadastrar Usuário
#      Cadastrar Usuário  -> Editar Usuário
#      (new)              -> Remover Usuário
# -----------------------------------------------------------------------
$listAnchor = $d.Content
$foundAnchor = $listAnchor.Find.Execute("Gerar Relatórios", $false, $false, $false, $false, $false, $true, 1, $false)
if (-not $foundAnchor) {
    throw "Could not find the 'Gerar Relatórios' list item."
}

$paraGerar = $listAnchor.Paragraphs(1)
$paraCadastrar = $paraGerar.Next()

# Rename "Gerar Relatórios" -> "Cadastrar Usuário" (keep the paragraph
# mark / list formatting untouched).
$rGerar = $d.Range($paraGerar.Range.Start, $paraGerar.Range.End - 1)
$rGerar.Text = "Cadastrar Usuário"

# Rename the existing "Cadastrar Usuário" item -> "Editar Usuário".
$rCadastrar = $d.Range($paraCadastrar.Range.Start, $paraCadastrar.Range.End - 1)
$rCadastrar.Text = "Editar Usuário"

# Add a new list item "Remover Usuário" right after it, inheriting the
# same paragraph (list) formatting.
$paraCadastrar.Range.InsertParagraphAfter()
$paraRemover = $paraCadastrar.Next()
$rRemover = $d.Range($paraRemover.Range.Start, $paraRemover.Range.End - 1)
$rRemover.Text = "Remover Usuário"
